$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new (blank) column at N ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns("N:N").Insert()
$wsRepay.Columns("N:N").ColumnWidth = 9.86

# --- Switch the active / selected sheet from "Transactions" to "Repayment schedule" ---
$wsRepay.Activate()
$wsRepay.Range("S8").Select()

# Transactions sheet keeps its own prior selection, just loses the "active tab" flag
# (handled automatically by activating the Repayment schedule sheet above).
